$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '48.350.78'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.31%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.523.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.49%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '109.75'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.90%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '322.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.40%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.553'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.86%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.44'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.40'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +11.82%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.83%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.27'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.61%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.919.27'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.531.94'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.853'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.99%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '48.188.08'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.30'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.77%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.64%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0948'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.94%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.12'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.65%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '274.43'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +11.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.81%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.99'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.41%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.59%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.14'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.87%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.144'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.84'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.59%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.67'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.75'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.85%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.55%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.44%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.67%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.56%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.02'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.38%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.52'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.23%  '

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.24'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.08%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.02'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.64%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0301'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.032.06'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.44%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.16'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.17%  '

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.80%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.89'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +7.40%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.31%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.15'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.52%  '
